$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("G3").Value = 2
$ws.Range("H3").Value = 3.1
$ws.Range("I3").Value = 4.33
$ws.Range("J3").Value = 2.75
$ws.Range("K3").Value = 1.95
$ws.Range("L3").Value = 4.75
$ws.Range("M3").Value = 1.1
$ws.Range("N3").Value = 7
$ws.Range("O3").Value = 1.5
$ws.Range("P3").Value = 2.5
$ws.Range("Q3").Value = 1.95
$ws.Range("R3").Value = 1.9
$ws.Range("S3").Value = 2.6
$ws.Range("T3").Value = 1.48
$ws.Range("U3").Value = 3.95
$ws.Range("V3").Value = 1.24
$ws.Range("W3").Value = 5
$ws.Range("X3").Value = 1.17
$ws.Range("Y3").Value = 1.57
$ws.Range("Z3").Value = 2.25
$ws.Range("AD3").Value = 8
$ws.Range("AE3").Value = 9.5
$ws.Range("AF3").Value = 17
$ws.Range("AG3").Value = 21
$ws.Range("AH3").Value = 41
$ws.Range("AI3").Value = 6.5
$ws.Range("AJ3").Value = 6
$ws.Range("AL3").Value = 81
$ws.Range("AN3").Value = 9
$ws.Range("AO3").Value = 19
$ws.Range("AP3").Value = 15
$ws.Range("AQ3").Value = 41
# Row 4
$ws.Range("G4").Value = 1.73
$ws.Range("H4").Value = 3.6
$ws.Range("I4").Value = 4.75
$ws.Range("J4").Value = 2.4
$ws.Range("K4").Value = 2.1
$ws.Range("S4").Value = 2.15
$ws.Range("T4").Value = 1.67
$ws.Range("U4").Value = 3.2
$ws.Range("V4").Value = 1.34
$ws.Range("Y4").Value = 1.44
$ws.Range("Z4").Value = 2.63
$ws.Range("AA4").Value = 2
$ws.Range("AB4").Value = 1.73
$ws.Range("AC4").Value = 6
$ws.Range("AD4").Value = 7.5
$ws.Range("AG4").Value = 15
$ws.Range("AI4").Value = 8.5
$ws.Range("AK4").Value = 19
$ws.Range("AP4").Value = 15
$ws.Range("AS4").Value = 41
# Row 5
$ws.Range("G5").Value = 4.25
$ws.Range("H5").Value = 3.25
$ws.Range("I5").Value = 1.82
$ws.Range("J5").Value = 4.6
$ws.Range("K5").Value = 2.05
$ws.Range("L5").Value = 2.42
$ws.Range("O5").Value = 1.25
$ws.Range("P5").Value = 3.25
$ws.Range("S5").Value = 1.75
$ws.Range("T5").Value = 1.87
$ws.Range("W5").Value = 2.7
$ws.Range("X5").Value = 1.35
$ws.Range("Y5").Value = 1.4
$ws.Range("Z5").Value = 2.55
$ws.Range("AA5").Value = 1.62
$ws.Range("AB5").Value = 2.02
$ws.Range("AG5").Value = 40
$ws.Range("AH5").Value = 40
$ws.Range("AI5").Value = 10.5
$ws.Range("AJ5").Value = 6.5
$ws.Range("AN5").Value = 7.7
$ws.Range("AO5").Value = 9.5
$ws.Range("AP5").Value = 7.9
$ws.Range("AQ5").Value = 16
$ws.Range("AR5").Value = 13.5
$ws.Range("AS5").Value = 22
